$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.306.66'
$ws.Range('E2').Value = '  +0.22%  '
$ws.Range('D3').Value = '1.858.60'
$ws.Range('E3').Value = '  -0.18%  '
$ws.Range('E4').Value = '  +0.28%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7017'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.12%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '238.28'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.60%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.001'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.20%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07878'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +2.95%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3031'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.49%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '24.49'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +4.98%  '
$ws.Range('E11').Value = '  +0.54%  '
$ws.Range('D12').Value = '1.866.34'
$ws.Range('E12').Value = '  +0.32%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.211'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.23%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.7074'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.39%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '89.55'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.08%  '
$ws.Range('D16').Value = '29.361.29'
$ws.Range('E16').Value = '  +0.35%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '5.792'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.83%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000007825'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.53%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.22'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.42%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '237.41'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.15%  '
$ws.Range('D21').Value = '2.126.81'
$ws.Range('E21').Value = '  +0.69%  '
$ws.Range('E22').Value = '  +0.08%  '
$ws.Range('E23').Value = '  +0.29%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '7.566'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.89%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '162.17'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.23%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.899'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.34%  '
$ws.Range('E27').Value = '  -2.79%  '
$ws.Range('E28').Value = '  +0.12%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.909'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -3.30%  '
$ws.Range('E30').Value = '  -0.69%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.479'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.11%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.303'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.13%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.034'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.01%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05188'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.01%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.180'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.75%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7109'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.47%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.9988'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.04%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.678'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.16%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01852'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.20%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.684'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.12%  '
$ws.Range('D41').Value = '1.144.11'
$ws.Range('E41').Value = '  +0.55%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.9216'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.66%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.952'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.25%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.4249'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.72%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '70.42'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.42%  '
$ws.Range('E46').Value = '  +0.29%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '102.67'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.66%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.5315'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.57%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.747'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -3.22%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.187'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.23%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.012'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.80%  '
